# Resuelto conflicto en db.sqlite y conservada la versión local
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the stock quantity for "Silla Gamer" (row 6, column D) from 3 to 0
$ws.Range("D6").Value = 0

# Restore the view state captured when the file was last saved:
# zoom level and the active cell/selection on the sheet
$excel.ActiveWindow.Zoom = 107
$ws.Range("G13").Select()
